# Guiglia.xlsx update: append daily COVID data rows through 2021-12-08 (commit: "aggiornamento fino a 8/12")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the workbook view settings (values unchanged by this edit).
$win = $excel.ActiveWindow
$win.AutoFilterDateGrouping = $true
$win.TabRatio = 600
$win.DisplayHorizontalScrollBar = $true
$win.DisplayVerticalScrollBar = $true

# Reuse the existing date-formatted style (column A, style index 2) for the new rows
# by copying formatting from the last existing data row (A385) onto the new range.
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new rows (386-464) with date serials (col A), new-positive counts (col B),
# 7-day rolling sum (col C) and the per-100k-inhabitants figure (col D).
$ws.Range("A386").Value2 = 44460
$ws.Range("B386").Value2 = 0
$ws.Range("C386").Value2 = 3
$ws.Range("D386").Value2 = 76.2970498474059
$ws.Range("A387").Value2 = 44461
$ws.Range("B387").Value2 = 0
$ws.Range("C387").Value2 = 3
$ws.Range("D387").Value2 = 76.2970498474059
$ws.Range("A388").Value2 = 44462
$ws.Range("B388").Value2 = 1
$ws.Range("C388").Value2 = 4
$ws.Range("D388").Value2 = 101.7293997965412
$ws.Range("A389").Value2 = 44463
$ws.Range("B389").Value2 = 3
$ws.Range("C389").Value2 = 5
$ws.Range("D389").Value2 = 127.1617497456765
$ws.Range("A390").Value2 = 44464
$ws.Range("B390").Value2 = 0
$ws.Range("C390").Value2 = 5
$ws.Range("D390").Value2 = 127.1617497456765
$ws.Range("A391").Value2 = 44465
$ws.Range("B391").Value2 = 0
$ws.Range("C391").Value2 = 5
$ws.Range("D391").Value2 = 127.1617497456765
$ws.Range("A392").Value2 = 44466
$ws.Range("B392").Value2 = 1
$ws.Range("C392").Value2 = 5
$ws.Range("D392").Value2 = 127.1617497456765
$ws.Range("A393").Value2 = 44467
$ws.Range("B393").Value2 = 0
$ws.Range("C393").Value2 = 5
$ws.Range("D393").Value2 = 127.1617497456765
$ws.Range("A394").Value2 = 44468
$ws.Range("B394").Value2 = 0
$ws.Range("C394").Value2 = 5
$ws.Range("D394").Value2 = 127.1617497456765
$ws.Range("A395").Value2 = 44469
$ws.Range("B395").Value2 = 2
$ws.Range("C395").Value2 = 6
$ws.Range("D395").Value2 = 152.5940996948118
$ws.Range("A396").Value2 = 44470
$ws.Range("B396").Value2 = 0
$ws.Range("C396").Value2 = 3
$ws.Range("D396").Value2 = 76.2970498474059
$ws.Range("A397").Value2 = 44471
$ws.Range("B397").Value2 = 0
$ws.Range("C397").Value2 = 3
$ws.Range("D397").Value2 = 76.2970498474059
$ws.Range("A398").Value2 = 44472
$ws.Range("B398").Value2 = 0
$ws.Range("C398").Value2 = 3
$ws.Range("D398").Value2 = 76.2970498474059
$ws.Range("A399").Value2 = 44473
$ws.Range("B399").Value2 = 0
$ws.Range("C399").Value2 = 2
$ws.Range("D399").Value2 = 50.8646998982706
$ws.Range("A400").Value2 = 44474
$ws.Range("B400").Value2 = 0
$ws.Range("C400").Value2 = 2
$ws.Range("D400").Value2 = 50.8646998982706
$ws.Range("A401").Value2 = 44475
$ws.Range("B401").Value2 = 0
$ws.Range("C401").Value2 = 2
$ws.Range("D401").Value2 = 50.8646998982706
$ws.Range("A402").Value2 = 44476
$ws.Range("B402").Value2 = 0
$ws.Range("C402").Value2 = 0
$ws.Range("D402").Value2 = 0
$ws.Range("A403").Value2 = 44477
$ws.Range("B403").Value2 = 0
$ws.Range("C403").Value2 = 0
$ws.Range("D403").Value2 = 0
$ws.Range("A404").Value2 = 44478
$ws.Range("B404").Value2 = 0
$ws.Range("C404").Value2 = 0
$ws.Range("D404").Value2 = 0
$ws.Range("A405").Value2 = 44479
$ws.Range("B405").Value2 = 0
$ws.Range("C405").Value2 = 0
$ws.Range("D405").Value2 = 0
$ws.Range("A406").Value2 = 44480
$ws.Range("B406").Value2 = 0
$ws.Range("C406").Value2 = 0
$ws.Range("D406").Value2 = 0
$ws.Range("A407").Value2 = 44481
$ws.Range("B407").Value2 = 0
$ws.Range("C407").Value2 = 0
$ws.Range("D407").Value2 = 0
$ws.Range("A408").Value2 = 44482
$ws.Range("B408").Value2 = 0
$ws.Range("C408").Value2 = 0
$ws.Range("D408").Value2 = 0
$ws.Range("A409").Value2 = 44483
$ws.Range("B409").Value2 = 0
$ws.Range("C409").Value2 = 0
$ws.Range("D409").Value2 = 0
$ws.Range("A410").Value2 = 44484
$ws.Range("B410").Value2 = 1
$ws.Range("C410").Value2 = 1
$ws.Range("D410").Value2 = 25.4323499491353
$ws.Range("A411").Value2 = 44485
$ws.Range("B411").Value2 = 1
$ws.Range("C411").Value2 = 2
$ws.Range("D411").Value2 = 50.8646998982706
$ws.Range("A412").Value2 = 44486
$ws.Range("B412").Value2 = 0
$ws.Range("C412").Value2 = 2
$ws.Range("D412").Value2 = 50.8646998982706
$ws.Range("A413").Value2 = 44487
$ws.Range("B413").Value2 = 1
$ws.Range("C413").Value2 = 3
$ws.Range("D413").Value2 = 76.2970498474059
$ws.Range("A414").Value2 = 44488
$ws.Range("B414").Value2 = 0
$ws.Range("C414").Value2 = 3
$ws.Range("D414").Value2 = 76.2970498474059
$ws.Range("A415").Value2 = 44489
$ws.Range("B415").Value2 = 0
$ws.Range("C415").Value2 = 3
$ws.Range("D415").Value2 = 76.2970498474059
$ws.Range("A416").Value2 = 44490
$ws.Range("B416").Value2 = 0
$ws.Range("C416").Value2 = 3
$ws.Range("D416").Value2 = 76.2970498474059
$ws.Range("A417").Value2 = 44491
$ws.Range("B417").Value2 = 0
$ws.Range("C417").Value2 = 2
$ws.Range("D417").Value2 = 50.8646998982706
$ws.Range("A418").Value2 = 44492
$ws.Range("B418").Value2 = 0
$ws.Range("C418").Value2 = 1
$ws.Range("D418").Value2 = 25.4323499491353
$ws.Range("A419").Value2 = 44493
$ws.Range("B419").Value2 = 0
$ws.Range("C419").Value2 = 1
$ws.Range("D419").Value2 = 25.4323499491353
$ws.Range("A420").Value2 = 44494
$ws.Range("B420").Value2 = 0
$ws.Range("C420").Value2 = 0
$ws.Range("D420").Value2 = 0
$ws.Range("A421").Value2 = 44495
$ws.Range("B421").Value2 = 0
$ws.Range("C421").Value2 = 0
$ws.Range("D421").Value2 = 0
$ws.Range("A422").Value2 = 44496
$ws.Range("B422").Value2 = 0
$ws.Range("C422").Value2 = 0
$ws.Range("D422").Value2 = 0
$ws.Range("A423").Value2 = 44497
$ws.Range("B423").Value2 = 0
$ws.Range("C423").Value2 = 0
$ws.Range("D423").Value2 = 0
$ws.Range("A424").Value2 = 44498
$ws.Range("B424").Value2 = 0
$ws.Range("C424").Value2 = 0
$ws.Range("D424").Value2 = 0
$ws.Range("A425").Value2 = 44499
$ws.Range("B425").Value2 = 0
$ws.Range("C425").Value2 = 0
$ws.Range("D425").Value2 = 0
$ws.Range("A426").Value2 = 44500
$ws.Range("B426").Value2 = 0
$ws.Range("C426").Value2 = 0
$ws.Range("D426").Value2 = 0
$ws.Range("A427").Value2 = 44501
$ws.Range("B427").Value2 = 0
$ws.Range("C427").Value2 = 0
$ws.Range("D427").Value2 = 0
$ws.Range("A428").Value2 = 44502
$ws.Range("B428").Value2 = 0
$ws.Range("C428").Value2 = 0
$ws.Range("D428").Value2 = 0
$ws.Range("A429").Value2 = 44503
$ws.Range("B429").Value2 = 0
$ws.Range("C429").Value2 = 0
$ws.Range("D429").Value2 = 0
$ws.Range("A430").Value2 = 44504
$ws.Range("B430").Value2 = 0
$ws.Range("C430").Value2 = 0
$ws.Range("D430").Value2 = 0
$ws.Range("A431").Value2 = 44505
$ws.Range("B431").Value2 = 0
$ws.Range("C431").Value2 = 0
$ws.Range("D431").Value2 = 0
$ws.Range("A432").Value2 = 44506
$ws.Range("B432").Value2 = 0
$ws.Range("C432").Value2 = 0
$ws.Range("D432").Value2 = 0
$ws.Range("A433").Value2 = 44507
$ws.Range("B433").Value2 = 0
$ws.Range("C433").Value2 = 0
$ws.Range("D433").Value2 = 0
$ws.Range("A434").Value2 = 44508
$ws.Range("B434").Value2 = 0
$ws.Range("C434").Value2 = 0
$ws.Range("D434").Value2 = 0
$ws.Range("A435").Value2 = 44509
$ws.Range("B435").Value2 = 0
$ws.Range("C435").Value2 = 0
$ws.Range("D435").Value2 = 0
$ws.Range("A436").Value2 = 44510
$ws.Range("B436").Value2 = 0
$ws.Range("C436").Value2 = 0
$ws.Range("D436").Value2 = 0
$ws.Range("A437").Value2 = 44511
$ws.Range("B437").Value2 = 0
$ws.Range("C437").Value2 = 0
$ws.Range("D437").Value2 = 0
$ws.Range("A438").Value2 = 44512
$ws.Range("B438").Value2 = 0
$ws.Range("C438").Value2 = 0
$ws.Range("D438").Value2 = 0
$ws.Range("A439").Value2 = 44513
$ws.Range("B439").Value2 = 0
$ws.Range("C439").Value2 = 0
$ws.Range("D439").Value2 = 0
$ws.Range("A440").Value2 = 44514
$ws.Range("B440").Value2 = 0
$ws.Range("C440").Value2 = 0
$ws.Range("D440").Value2 = 0
$ws.Range("A441").Value2 = 44515
$ws.Range("B441").Value2 = 2
$ws.Range("C441").Value2 = 2
$ws.Range("D441").Value2 = 50.8646998982706
$ws.Range("A442").Value2 = 44516
$ws.Range("B442").Value2 = 4
$ws.Range("C442").Value2 = 6
$ws.Range("D442").Value2 = 152.5940996948118
$ws.Range("A443").Value2 = 44517
$ws.Range("B443").Value2 = 0
$ws.Range("C443").Value2 = 6
$ws.Range("D443").Value2 = 152.5940996948118
$ws.Range("A444").Value2 = 44518
$ws.Range("B444").Value2 = 1
$ws.Range("C444").Value2 = 7
$ws.Range("D444").Value2 = 178.0264496439471
$ws.Range("A445").Value2 = 44519
$ws.Range("B445").Value2 = 0
$ws.Range("C445").Value2 = 7
$ws.Range("D445").Value2 = 178.0264496439471
$ws.Range("A446").Value2 = 44520
$ws.Range("B446").Value2 = 0
$ws.Range("C446").Value2 = 7
$ws.Range("D446").Value2 = 178.0264496439471
$ws.Range("A447").Value2 = 44521
$ws.Range("B447").Value2 = 0
$ws.Range("C447").Value2 = 7
$ws.Range("D447").Value2 = 178.0264496439471
$ws.Range("A448").Value2 = 44522
$ws.Range("B448").Value2 = 0
$ws.Range("C448").Value2 = 5
$ws.Range("D448").Value2 = 127.1617497456765
$ws.Range("A449").Value2 = 44523
$ws.Range("B449").Value2 = 0
$ws.Range("C449").Value2 = 1
$ws.Range("D449").Value2 = 25.4323499491353
$ws.Range("A450").Value2 = 44524
$ws.Range("B450").Value2 = 2
$ws.Range("C450").Value2 = 3
$ws.Range("D450").Value2 = 76.2970498474059
$ws.Range("A451").Value2 = 44525
$ws.Range("B451").Value2 = 0
$ws.Range("C451").Value2 = 2
$ws.Range("D451").Value2 = 50.8646998982706
$ws.Range("A452").Value2 = 44526
$ws.Range("B452").Value2 = 1
$ws.Range("C452").Value2 = 3
$ws.Range("D452").Value2 = 76.2970498474059
$ws.Range("A453").Value2 = 44527
$ws.Range("B453").Value2 = 0
$ws.Range("C453").Value2 = 3
$ws.Range("D453").Value2 = 76.2970498474059
$ws.Range("A454").Value2 = 44528
$ws.Range("B454").Value2 = 0
$ws.Range("C454").Value2 = 3
$ws.Range("D454").Value2 = 76.2970498474059
$ws.Range("A455").Value2 = 44529
$ws.Range("B455").Value2 = 0
$ws.Range("C455").Value2 = 3
$ws.Range("D455").Value2 = 76.2970498474059
$ws.Range("A456").Value2 = 44530
$ws.Range("B456").Value2 = 1
$ws.Range("C456").Value2 = 4
$ws.Range("D456").Value2 = 101.7293997965412
$ws.Range("A457").Value2 = 44531
$ws.Range("B457").Value2 = 0
$ws.Range("C457").Value2 = 2
$ws.Range("D457").Value2 = 50.8646998982706
$ws.Range("A458").Value2 = 44532
$ws.Range("B458").Value2 = 0
$ws.Range("C458").Value2 = 2
$ws.Range("D458").Value2 = 50.8646998982706
$ws.Range("A459").Value2 = 44533
$ws.Range("B459").Value2 = 0
$ws.Range("C459").Value2 = 1
$ws.Range("D459").Value2 = 25.4323499491353
$ws.Range("A460").Value2 = 44534
$ws.Range("B460").Value2 = 0
$ws.Range("C460").Value2 = 1
$ws.Range("D460").Value2 = 25.4323499491353
$ws.Range("A461").Value2 = 44535
$ws.Range("B461").Value2 = 5
$ws.Range("C461").Value2 = 6
$ws.Range("D461").Value2 = 152.5940996948118
$ws.Range("A462").Value2 = 44536
$ws.Range("B462").Value2 = 1
$ws.Range("C462").Value2 = 7
$ws.Range("D462").Value2 = 178.0264496439471
$ws.Range("A463").Value2 = 44537
$ws.Range("B463").Value2 = 5
$ws.Range("C463").Value2 = 11
$ws.Range("D463").Value2 = 279.7558494404883
$ws.Range("A464").Value2 = 44538
$ws.Range("B464").Value2 = 1
$ws.Range("C464").Value2 = 12
$ws.Range("D464").Value2 = 305.1881993896236
